$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.611780881881714
$ws.Range("B1").Value = 4.168026924133301
$ws.Range("C1").Value = 3.657693386077881
$ws.Range("D1").Value = 1.80646812915802
$ws.Range("E1").Value = 1.03722095489502
